$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $oldStyle = $r.Style
    $r.NumberFormat = '@'
    $r.Value = $val
    $r.Style = $oldStyle
}

$ws.Range('D2').Value = '64.195.48'
$ws.Range('E2').Value = '  -3.89%  '
$ws.Range('D3').Value = '3.159.40'
$ws.Range('E3').Value = '  -8.78%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '561.76'
$ws.Range('E5').Value = '  -4.15%  '
Set-TextValue 'D6' '168.43'
$ws.Range('E6').Value = '  -5.68%  '
Set-TextValue 'D7' '0.609'
$ws.Range('E7').Value = '  -3.52%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '3.159.88'
$ws.Range('E9').Value = '  -8.84%  '
$ws.Range('E10').Value = '  -7.60%  '
$ws.Range('E11').Value = '  -5.83%  '
$ws.Range('E12').Value = '  -6.23%  '
$ws.Range('D13').Value = '3.701.37'
$ws.Range('E13').Value = '  -8.99%  '
$ws.Range('E14').Value = '  +0.60%  '
Set-TextValue 'D15' '27.11'
$ws.Range('E15').Value = '  -10.32%  '
$ws.Range('D16').Value = '64.169.50'
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('E17').Value = '  -6.61%  '
$ws.Range('D18').Value = '3.153.52'
$ws.Range('E18').Value = '  -10.24%  '
$ws.Range('E19').Value = '  -4.32%  '
Set-TextValue 'D20' '12.90'
$ws.Range('E20').Value = '  -7.29%  '
Set-TextValue 'D21' '351.77'
$ws.Range('E21').Value = '  -5.68%  '
Set-TextValue 'D22' '7.19'
$ws.Range('E22').Value = '  -6.38%  '
$ws.Range('E23').Value = '  +0.44%  '
Set-TextValue 'D24' '68.07'
$ws.Range('E24').Value = '  -7.52%  '
Set-TextValue 'D25' '0.498'
$ws.Range('E25').Value = '  -7.08%  '
$ws.Range('E26').Value = '  -10.26%  '
Set-TextValue 'D27' '9.54'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('E29').Value = '  +0.04%  '
Set-TextValue 'D30' '0.999'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  -6.11%  '
Set-TextValue 'D32' '5.45'
$ws.Range('E32').Value = '  -8.04%  '
Set-TextValue 'D33' '21.89'
$ws.Range('E33').Value = '  -7.70%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D34' '1.20'
$ws.Range('E34').Value = '  -6.61%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D35' '6.59'
$ws.Range('E35').Value = '  -7.19%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D36' '154.38'
$ws.Range('E36').Value = '  -5.19%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.42'
$ws.Range('E37').Value = '  -9.59%  '
Set-TextValue 'D38' '0.816'
$ws.Range('E38').Value = '  -7.87%  '
Set-TextValue 'D39' '26.43'
$ws.Range('E39').Value = '  -5.75%  '
$ws.Range('E40').Value = '  -7.51%  '
$ws.Range('D41').Value = '2.638.73'
$ws.Range('E41').Value = '  -4.58%  '
Set-TextValue 'D42' '2.45'
$ws.Range('E42').Value = '  -5.94%  '
$ws.Range('E43').Value = '  -7.99%  '
Set-TextValue 'D44' '39.36'
$ws.Range('E44').Value = '  -1.83%  '
Set-TextValue 'D45' '5.98'
$ws.Range('E45').Value = '  -7.86%  '
Set-TextValue 'D46' '0.0648'
$ws.Range('E46').Value = '  -7.09%  '
Set-TextValue 'D47' '23.61'
$ws.Range('E47').Value = '  -7.20%  '
Set-TextValue 'D48' '320.41'
$ws.Range('E48').Value = '  -6.21%  '
Set-TextValue 'D49' '0.0270'
$ws.Range('E49').Value = '  -6.60%  '
$ws.Range('E50').Value = '  -3.64%  '
Set-TextValue 'D51' '0.999'
$ws.Range('E51').Value = '  -0.01%  '
